$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(5, 8).Value2 = 367.5
$ws.Cells.Item(5, 9).Value2 = 528
$ws.Cells.Item(5, 11).Value2 = 528
$ws.Cells.Item(5, 13).Value2 = -413
$ws.Cells.Item(28, 8).Value2 = 460.03125
$ws.Cells.Item(28, 9).Value2 = 442.86667
$ws.Cells.Item(28, 10).Value2 = 475.17648
$ws.Cells.Item(28, 11).Value2 = 442.86667
$ws.Cells.Item(28, 12).Value2 = 475.17648
$ws.Cells.Item(28, 13).Value2 = 42.13333
$ws.Cells.Item(28, 14).Value2 = -1445.17648
$ws.Cells.Item(40, 8).Value2 = 5368.885
$ws.Cells.Item(40, 9).Value2 = 5492.44
$ws.Cells.Item(40, 10).Value2 = 2280
$ws.Cells.Item(40, 11).Value2 = 5492.44
$ws.Cells.Item(40, 12).Value2 = 2280
$ws.Cells.Item(40, 13).Value2 = -5317.44
$ws.Cells.Item(40, 14).Value2 = -2630
$ws.Cells.Item(64, 8).Value2 = 3004.4314
$ws.Cells.Item(64, 9).Value2 = 3014.2856
$ws.Cells.Item(64, 10).Value2 = 2992.4348
$ws.Cells.Item(64, 11).Value2 = 3014.2856
$ws.Cells.Item(64, 12).Value2 = 2992.4348
$ws.Cells.Item(64, 13).Value2 = -2766.2856
$ws.Cells.Item(64, 14).Value2 = -3488.4348
$ws.Cells.Item(67, 8).Value2 = 3004.4314
$ws.Cells.Item(67, 9).Value2 = 3014.2856
$ws.Cells.Item(67, 10).Value2 = 2992.4348
$ws.Cells.Item(67, 11).Value2 = 3014.2856
$ws.Cells.Item(67, 12).Value2 = 2992.4348
$ws.Cells.Item(67, 13).Value2 = -2156.2856
$ws.Cells.Item(67, 14).Value2 = -4708.4348
$ws.Cells.Item(76, 8).Value2 = 3104.4783
$ws.Cells.Item(76, 9).Value2 = 3095.3809
$ws.Cells.Item(76, 11).Value2 = 3095.3809
$ws.Cells.Item(76, 13).Value2 = -2780.3809
$ws.Cells.Item(79, 8).Value2 = 3104.4783
$ws.Cells.Item(79, 9).Value2 = 3095.3809
$ws.Cells.Item(79, 11).Value2 = 3095.3809
$ws.Cells.Item(79, 13).Value2 = -2003.3809
$ws.Cells.Item(107, 8).Value2 = 5087.25
$ws.Cells.Item(107, 9).Value2 = 6366.8887
$ws.Cells.Item(107, 10).Value2 = 1248.3334
$ws.Cells.Item(107, 11).Value2 = 6366.8887
$ws.Cells.Item(107, 12).Value2 = 1248.3334
$ws.Cells.Item(107, 13).Value2 = -4446.8887
$ws.Cells.Item(107, 14).Value2 = -5088.3334
$ws.Cells.Item(125, 8).Value2 = 1567.5714
$ws.Cells.Item(125, 9).Value2 = 1349
$ws.Cells.Item(125, 10).Value2 = 1859
$ws.Cells.Item(125, 11).Value2 = 12141
$ws.Cells.Item(125, 12).Value2 = 16731
$ws.Cells.Item(125, 13).Value2 = -9681
$ws.Cells.Item(125, 14).Value2 = -21651
$ws.Cells.Item(126, 8).Value2 = 48000
$ws.Cells.Item(126, 9).Value2 = 0
$ws.Cells.Item(126, 10).Value2 = 48000
$ws.Cells.Item(126, 11).Value2 = 0
$ws.Cells.Item(126, 12).Value2 = 48000
$ws.Cells.Item(126, 14).Value2 = -57880
$ws.Cells.Item(127, 8).Value2 = 1119.4286
$ws.Cells.Item(127, 9).Value2 = 940.0909
$ws.Cells.Item(127, 10).Value2 = 1777
$ws.Cells.Item(127, 11).Value2 = 2820.2727
$ws.Cells.Item(127, 12).Value2 = 5331
$ws.Cells.Item(127, 13).Value2 = 2139.7273
$ws.Cells.Item(127, 14).Value2 = -15251
$ws.Cells.Item(128, 8).Value2 = 58774.332
$ws.Cells.Item(128, 9).Value2 = 0
$ws.Cells.Item(128, 10).Value2 = 58774.332
$ws.Cells.Item(128, 11).Value2 = 0
$ws.Cells.Item(128, 12).Value2 = 58774.332
$ws.Cells.Item(128, 14).Value2 = -68734.33199999999
$ws.Cells.Item(129, 8).Value2 = 1099.9048
$ws.Cells.Item(129, 9).Value2 = 1206.1538
$ws.Cells.Item(129, 10).Value2 = 1052.2759
$ws.Cells.Item(129, 11).Value2 = 3618.4614
$ws.Cells.Item(129, 12).Value2 = 3156.8277
$ws.Cells.Item(129, 13).Value2 = 1381.5386
$ws.Cells.Item(129, 14).Value2 = -13156.8277
$ws.Cells.Item(130, 8).Value2 = 0
$ws.Cells.Item(130, 9).Value2 = 0
$ws.Cells.Item(130, 10).Value2 = 0
$ws.Cells.Item(130, 11).Value2 = 0
$ws.Cells.Item(130, 12).Value2 = 0
$ws.Cells.Item(131, 8).Value2 = 3575
$ws.Cells.Item(131, 9).Value2 = 3990
$ws.Cells.Item(131, 10).Value2 = 3461.818
$ws.Cells.Item(131, 11).Value2 = 11970
$ws.Cells.Item(131, 12).Value2 = 10385.454
$ws.Cells.Item(131, 13).Value2 = -6930
$ws.Cells.Item(131, 14).Value2 = -20465.454
$ws.Cells.Item(132, 8).Value2 = 14119.889
$ws.Cells.Item(132, 9).Value2 = 1867.4445
$ws.Cells.Item(132, 10).Value2 = 99887
$ws.Cells.Item(132, 11).Value2 = 5602.333500000001
$ws.Cells.Item(132, 12).Value2 = 299661
$ws.Cells.Item(132, 13).Value2 = -3072.333500000001
$ws.Cells.Item(132, 14).Value2 = -304721
$ws.Cells.Item(133, 8).Value2 = 49738.57
$ws.Cells.Item(133, 9).Value2 = 0
$ws.Cells.Item(133, 10).Value2 = 49738.57
$ws.Cells.Item(133, 11).Value2 = 0
$ws.Cells.Item(133, 12).Value2 = 49738.57
$ws.Cells.Item(133, 14).Value2 = -59858.57
$ws.Cells.Item(134, 8).Value2 = 41600
$ws.Cells.Item(134, 9).Value2 = 0
$ws.Cells.Item(134, 10).Value2 = 41600
$ws.Cells.Item(134, 11).Value2 = 0
$ws.Cells.Item(134, 12).Value2 = 41600
$ws.Cells.Item(134, 14).Value2 = -51740
$ws.Cells.Item(135, 8).Value2 = 17858168
$ws.Cells.Item(135, 9).Value2 = 1070.7084
$ws.Cells.Item(135, 10).Value2 = 125000744
$ws.Cells.Item(135, 11).Value2 = 9636.375599999999
$ws.Cells.Item(135, 12).Value2 = 1125006696
$ws.Cells.Item(135, 13).Value2 = -7101.375599999999
$ws.Cells.Item(135, 14).Value2 = -1125011766
$ws.Cells.Item(136, 8).Value2 = 45522.223
$ws.Cells.Item(136, 9).Value2 = 0
$ws.Cells.Item(136, 10).Value2 = 45522.223
$ws.Cells.Item(136, 11).Value2 = 0
$ws.Cells.Item(136, 12).Value2 = 45522.223
$ws.Cells.Item(136, 14).Value2 = -55722.223
$ws.Cells.Item(137, 8).Value2 = 1883180.9
$ws.Cells.Item(137, 9).Value2 = 3503491.2
$ws.Cells.Item(137, 10).Value2 = 7032
$ws.Cells.Item(137, 11).Value2 = 10510473.6
$ws.Cells.Item(137, 12).Value2 = 21096
$ws.Cells.Item(137, 13).Value2 = -10507923.6
$ws.Cells.Item(137, 14).Value2 = -26196
$ws.Cells.Item(138, 8).Value2 = 2890.76
$ws.Cells.Item(138, 9).Value2 = 2004.3334
$ws.Cells.Item(138, 10).Value2 = 3532.6553
$ws.Cells.Item(138, 11).Value2 = 6013.0002
$ws.Cells.Item(138, 12).Value2 = 10597.9659
$ws.Cells.Item(138, 13).Value2 = -873.0002000000004
$ws.Cells.Item(138, 14).Value2 = -20877.9659
$ws.Cells.Item(139, 8).Value2 = 58333.332
$ws.Cells.Item(139, 9).Value2 = 0
$ws.Cells.Item(139, 10).Value2 = 58333.332
$ws.Cells.Item(139, 11).Value2 = 0
$ws.Cells.Item(139, 12).Value2 = 58333.332
$ws.Cells.Item(139, 14).Value2 = -68613.33199999999
$ws.Cells.Item(140, 8).Value2 = 38723.2
$ws.Cells.Item(140, 9).Value2 = 0
$ws.Cells.Item(140, 10).Value2 = 38723.2
$ws.Cells.Item(140, 11).Value2 = 0
$ws.Cells.Item(140, 12).Value2 = 38723.2
$ws.Cells.Item(140, 14).Value2 = -49083.2
$ws.Cells.Item(141, 8).Value2 = 3844.923
$ws.Cells.Item(141, 9).Value2 = 2267
$ws.Cells.Item(141, 10).Value2 = 9104.666999999999
$ws.Cells.Item(141, 11).Value2 = 6801
$ws.Cells.Item(141, 12).Value2 = 27314.001
$ws.Cells.Item(141, 13).Value2 = -1621
$ws.Cells.Item(141, 14).Value2 = -37674.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 10096.869
$ws.Cells.Item(32, 9).Value2 = 9697.429
$ws.Cells.Item(32, 10).Value2 = 11817.538
$ws.Cells.Item(32, 11).Value2 = 9697.429
$ws.Cells.Item(32, 12).Value2 = 11817.538
$ws.Cells.Item(32, 13).Value2 = -9410.429
$ws.Cells.Item(32, 14).Value2 = -12391.538
$ws.Cells.Item(63, 8).Value2 = 3611
$ws.Cells.Item(63, 10).Value2 = 5748
$ws.Cells.Item(63, 12).Value2 = 5748
$ws.Cells.Item(63, 14).Value2 = -7120
$ws.Cells.Item(66, 8).Value2 = 3611
$ws.Cells.Item(66, 10).Value2 = 5748
$ws.Cells.Item(66, 12).Value2 = 28740
$ws.Cells.Item(66, 14).Value2 = -35604
$ws.Cells.Item(132, 8).Value2 = 7694176.5
$ws.Cells.Item(132, 9).Value2 = 10001768
$ws.Cells.Item(132, 10).Value2 = 2205.7334
$ws.Cells.Item(132, 11).Value2 = 30005304
$ws.Cells.Item(132, 12).Value2 = 6617.2002
$ws.Cells.Item(132, 13).Value2 = -30002774
$ws.Cells.Item(132, 14).Value2 = -11677.2002

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value2 = 1380.2
$ws.Cells.Item(22, 9).Value2 = 1380.2
$ws.Cells.Item(22, 11).Value2 = 1380.2
$ws.Cells.Item(22, 13).Value2 = -1207.2
$ws.Cells.Item(86, 8).Value2 = 2552.2354
$ws.Cells.Item(86, 9).Value2 = 2164.0908
$ws.Cells.Item(86, 10).Value2 = 3263.8333
$ws.Cells.Item(86, 11).Value2 = 2164.0908
$ws.Cells.Item(86, 12).Value2 = 3263.8333
$ws.Cells.Item(86, 13).Value2 = -1041.0908
$ws.Cells.Item(86, 14).Value2 = -5509.8333
$ws.Cells.Item(89, 8).Value2 = 2552.2354
$ws.Cells.Item(89, 9).Value2 = 2164.0908
$ws.Cells.Item(89, 10).Value2 = 3263.8333
$ws.Cells.Item(89, 11).Value2 = 10820.454
$ws.Cells.Item(89, 12).Value2 = 16319.1665
$ws.Cells.Item(89, 13).Value2 = -5204.454
$ws.Cells.Item(89, 14).Value2 = -27551.1665
$ws.Cells.Item(134, 8).Value2 = 2278.5952
$ws.Cells.Item(134, 9).Value2 = 1823.9395
$ws.Cells.Item(134, 10).Value2 = 3945.6667
$ws.Cells.Item(134, 11).Value2 = 5471.818499999999
$ws.Cells.Item(134, 12).Value2 = 11837.0001
$ws.Cells.Item(134, 13).Value2 = -2936.818499999999
$ws.Cells.Item(134, 14).Value2 = -16907.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value2 = 1719.4
$ws.Cells.Item(22, 9).Value2 = 372.8
$ws.Cells.Item(22, 10).Value2 = 2392.7
$ws.Cells.Item(22, 11).Value2 = 372.8
$ws.Cells.Item(22, 12).Value2 = 2392.7
$ws.Cells.Item(22, 13).Value2 = -22.80000000000001
$ws.Cells.Item(22, 14).Value2 = -3092.7
$ws.Cells.Item(31, 8).Value2 = 3440046.8
$ws.Cells.Item(31, 9).Value2 = 2253.7742
$ws.Cells.Item(31, 10).Value2 = 5054767.5
$ws.Cells.Item(31, 11).Value2 = 2253.7742
$ws.Cells.Item(31, 12).Value2 = 5054767.5
$ws.Cells.Item(31, 13).Value2 = -1958.7742
$ws.Cells.Item(31, 14).Value2 = -5055357.5
$ws.Cells.Item(34, 8).Value2 = 3440046.8
$ws.Cells.Item(34, 9).Value2 = 2253.7742
$ws.Cells.Item(34, 10).Value2 = 5054767.5
$ws.Cells.Item(34, 11).Value2 = 2253.7742
$ws.Cells.Item(34, 12).Value2 = 5054767.5
$ws.Cells.Item(34, 13).Value2 = -2051.7742
$ws.Cells.Item(34, 14).Value2 = -5055171.5
$ws.Cells.Item(62, 8).Value2 = 3387.5833
$ws.Cells.Item(62, 9).Value2 = 3405.6667
$ws.Cells.Item(62, 10).Value2 = 3333.3333
$ws.Cells.Item(62, 11).Value2 = 3405.6667
$ws.Cells.Item(62, 12).Value2 = 3333.3333
$ws.Cells.Item(62, 13).Value2 = -2781.6667
$ws.Cells.Item(62, 14).Value2 = -4581.3333
$ws.Cells.Item(65, 8).Value2 = 3387.5833
$ws.Cells.Item(65, 9).Value2 = 3405.6667
$ws.Cells.Item(65, 10).Value2 = 3333.3333
$ws.Cells.Item(65, 11).Value2 = 17028.3335
$ws.Cells.Item(65, 12).Value2 = 16666.6665
$ws.Cells.Item(65, 13).Value2 = -13908.3335
$ws.Cells.Item(65, 14).Value2 = -22906.6665
$ws.Cells.Item(132, 8).Value2 = 335650.25
$ws.Cells.Item(132, 9).Value2 = 1605.9032
$ws.Cells.Item(132, 10).Value2 = 1277047.9
$ws.Cells.Item(132, 11).Value2 = 4817.7096
$ws.Cells.Item(132, 12).Value2 = 3831143.7
$ws.Cells.Item(132, 13).Value2 = -2287.7096
$ws.Cells.Item(132, 14).Value2 = -3836203.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(68, 8).Value2 = 1614.84
$ws.Cells.Item(68, 9).Value2 = 1117.08
$ws.Cells.Item(68, 10).Value2 = 1780.76
$ws.Cells.Item(68, 11).Value2 = 3351.24
$ws.Cells.Item(68, 12).Value2 = 5342.28
$ws.Cells.Item(68, 13).Value2 = -2540.24
$ws.Cells.Item(68, 14).Value2 = -6964.28
$ws.Cells.Item(71, 8).Value2 = 1614.84
$ws.Cells.Item(71, 9).Value2 = 1117.08
$ws.Cells.Item(71, 10).Value2 = 1780.76
$ws.Cells.Item(71, 11).Value2 = 10053.72
$ws.Cells.Item(71, 12).Value2 = 16026.84
$ws.Cells.Item(71, 13).Value2 = -5997.719999999999
$ws.Cells.Item(71, 14).Value2 = -24138.84
$ws.Cells.Item(113, 8).Value2 = 212900.48
$ws.Cells.Item(113, 9).Value2 = 280108.16
$ws.Cells.Item(113, 10).Value2 = 937.8461
$ws.Cells.Item(113, 11).Value2 = 840324.48
$ws.Cells.Item(113, 12).Value2 = 2813.5383
$ws.Cells.Item(113, 13).Value2 = -838154.48
$ws.Cells.Item(113, 14).Value2 = -7153.5383

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value2 = 708748.9399999999
$ws.Cells.Item(2, 9).Value2 = 860609.4399999999
$ws.Cells.Item(2, 10).Value2 = 66.666664
$ws.Cells.Item(2, 11).Value2 = 860609.4399999999
$ws.Cells.Item(2, 12).Value2 = 66.666664
$ws.Cells.Item(2, 13).Value2 = -860496.4399999999
$ws.Cells.Item(2, 14).Value2 = -292.666664
$ws.Cells.Item(132, 8).Value2 = 13335764
$ws.Cells.Item(132, 9).Value2 = 15386888
$ws.Cells.Item(132, 10).Value2 = 3460.3
$ws.Cells.Item(132, 11).Value2 = 46160664
$ws.Cells.Item(132, 12).Value2 = 10380.9
$ws.Cells.Item(132, 13).Value2 = -46158134
$ws.Cells.Item(132, 14).Value2 = -15440.9

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value2 = 726115.5600000001
$ws.Cells.Item(132, 9).Value2 = 989222.8
$ws.Cells.Item(132, 10).Value2 = 2570.6875
$ws.Cells.Item(132, 11).Value2 = 2967668.4
$ws.Cells.Item(132, 12).Value2 = 7712.0625
$ws.Cells.Item(132, 13).Value2 = -2965138.4
$ws.Cells.Item(132, 14).Value2 = -12772.0625
